$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the reference-number citations for the four authors whose
# bracketed numbers shifted (figures renumbered for references).
$ws.Range("B6").Value  = "Holobinko [75]"
$ws.Range("B7").Value  = "Juarez [76]"
$ws.Range("B8").Value  = "Lustig [77]"
$ws.Range("B3").Value  = "Engel [74]"

# Move the active selection to B3, matching the author's last edit location.
$ws.Range("B3").Select()
